# Applies the textual edits described by the diff to the active document.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "October  29, 2021 (07:12:55 PM)" "November   1, 2021 (05:06:20 PM)"

Replace-Text "First few lines are given as examples, your task is to complete the rest of the table." "The first few lines are given as examples, your task is to complete the rest of the table."

Replace-Text "If the provided input is not valid, request new input from user until user provides valid input." "If the provided input is not valid, request new input from the user until the user provides valid input."

Replace-Text "then check the user’s input. Consider any of these variations to mean yes:" "then checks the user’s input. Consider any of these variations to mean yes:"

Replace-Text ". Once user enters yes, exit the loop." ". Once the user enters yes, exit the loop."

Replace-Text "Ask the user to enter a positive integer, between 2 and 100 (including these boundary values 2 and 100). Make sure the value user enters is between these bounds. Then compute the sum of integers starting from 1 up to the integer user entered, and finally display that sum. Here are examples:" "Ask the user to enter a positive integer, between 2 and 100 (including the boundary values 2 and 100). Make sure the value the user enters is between these bounds. Then compute the sum of integers starting from 1 up to the integer user entered, and finally display that sum. Here are examples:"

Replace-Text "if user enters 5, compute: 1 + 2 + 3 + 4 + 5, then display 15 at the screen" "if the user enters 5, compute: 1 + 2 + 3 + 4 + 5, then display 15 on the screen"

Replace-Text "if user enters 8, compute: 1 + 2 + 3 + 4 + 5 + 6 + 7 + 8, then display 36 at the screen" "if the user enters 8, compute: 1 + 2 + 3 + 4 + 5 + 6 + 7 + 8, then display 36 on the screen"

Replace-Text "Do this next problem using" "Do this next problem using the"

Replace-Text "type. Ask the user to enter any numbers (can be positive, negative, or zero). Ignore all non-numeric inputs. Choose an appropriate sentinel value to enable user to indicate when they are done. Compute and display the average of all numbers user entered. If user entered no numbers, display" "type. Ask the user to enter any numbers (can be positive, negative, or zero). Ignore all non-numeric inputs. Choose an appropriate sentinel value to enable the user to indicate when they are done. Compute and display the average of all the numbers the user entered. If the user entered no numbers, display"

Replace-Text "solution, that mixes classes and decision structures." "solution, which mixes classes and decision structures."

Replace-Text "to add the following validation features:" "solution to add the following validation features:"

Write-Host "Done."
